$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Лист1")

$ws1.Range("C2").Value = 11
$ws1.Range("D2").Value = 10
$ws1.Range("E2").Value = 4

$ws1.Range("A3").Value = 9
$ws1.Range("C3").Value = 5
$ws1.Range("D3").Value = 3
$ws1.Range("E3").Value = 4

$ws1.Range("A4").Value = 6
$ws1.Range("B4").Value = 9
$ws1.Range("C4").Value = 14
$ws1.Range("D4").Value = 0
$ws1.Range("E4").Value = 0

$ws1.Range("A5").Value = 8
$ws1.Range("B5").Value = 5
$ws1.Range("C5").Value = 4
$ws1.Range("D5").Value = 0
$ws1.Range("E5").Value = 0

$ws1.Range("A6").Value = 12
$ws1.Range("D6").Value = 0
$ws1.Range("E6").Value = 0

$ws1.Range("D5").Interior.Color = 65535
$ws1.Range("E6").Interior.Color = 65535

$ws1.Range("G9").Select()
